$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for "Poroto granado" that belongs
# between the existing row 19 (2021-01-08) and the old row 20 (2022-02-04).
# Insert a fresh row at position 20 - this shifts the old rows 20..30 down
# to 21..31, matching the diff exactly - then populate the new row.
$ws.Rows(20).Insert()

$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("D20").Value = 44603
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 100112030
$ws.Range("G20").Value = "Poroto granado"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 80
$ws.Range("K20").Value = 33000
$ws.Range("L20").Value = 33000
$ws.Range("M20").Value = 33000
$ws.Range("N20").Value = "$/saco 25 kilos"
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 1320
$ws.Range("Q20").Value = 25
$ws.Range("R20").Value = "Hortaliza"
